# Update "summary of paper-excel version.xlsx"
# - Sheet1: split the combined "Title: ..." cell into a Title/sub-title row
#   plus a new citation row, pushing the existing Q&A rows down by one.
# - Sheet2: no longer the active tab; Sheet1 becomes active instead.
# - Defined name _GoBack now points at Sheet1!$A$4 (was $A$3).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1. Sheet1: insert a blank row under the title row, then populate row1/row2.
# ---------------------------------------------------------------------------
$ws1.Rows("2:2").Insert()

# Remove the old A1:B1 merge (the title is no longer a single merged cell).
$ws1.Range("A1:B1").UnMerge()

# Row 1: "Title:" label (bold) in A1, the paper title (regular) in B1.
$ws1.Range("A1").Value2 = "Title:"
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("A1").HorizontalAlignment = -4131
$ws1.Range("A1").IndentLevel = 5

$ws1.Range("B1").Value2 = " Improved performance of serially connected Li-ion batteries with active cell balancing in electric vehicles."
$ws1.Range("B1").Font.Bold = $false
$ws1.Range("B1").VerticalAlignment = -4160
$ws1.Range("B1").WrapText = $true

# Row 2: blank label cell (keeps the bold header look), citation in B2.
$ws1.Range("A2").Font.Bold = $true
$ws1.Range("A2").HorizontalAlignment = -4131
$ws1.Range("A2").IndentLevel = 5

$ws1.Range("B2").Value2 = "Einhorn, M., Roessler, W., & Fleig, J. (2011). Improved performance of serially connected Li-ion batteries with active cell balancing in electric vehicles. IEEE Transactions on Vehicular Technology, 60(6), 2448-2457."
$ws1.Range("B2").WrapText = $true
$ws1.Range("B2").Font.Name = "Arial"
$ws1.Range("B2").Font.Size = 10
$ws1.Range("B2").Font.Color = 2236962

$cite = $ws1.Range("B2")
$cite.Characters(1, 97).Font.Italic = $false
$cite.Characters(98, 34).Font.Italic = $true
$cite.Characters(132, 2).Font.Italic = $false
$cite.Characters(134, 2).Font.Italic = $true
$cite.Characters(136, 14).Font.Italic = $false

# ---------------------------------------------------------------------------
# 2. Column layout: A/B narrower, a new column C added (matches new widths).
# ---------------------------------------------------------------------------
$ws1.Columns("A").ColumnWidth = 86.43
$ws1.Columns("B").ColumnWidth = 54.57
$ws1.Columns("C").ColumnWidth = 62.71

# ---------------------------------------------------------------------------
# 3. Defined name _GoBack now refers to Sheet1!$A$4.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -match "_GoBack") {
        $n.RefersTo = "=Sheet1!`$A`$4"
    }
}

# ---------------------------------------------------------------------------
# 4. Active sheet / selections: Sheet1 becomes active (tabSelected), Sheet2
#    keeps a B2:B16 selection but is no longer the active tab.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B2:B16").Select()

$ws1.Activate()
$ws1.Range("A20").Select()
